$d = $word.ActiveDocument

# Find the Subtitle paragraph ("Laws, policies, and practical realities") so the
# two new Author paragraphs get inserted right after it (and before the Date
# paragraph), regardless of its exact index.
$subtitle = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Style.NameLocal -eq "Subtitle") {
        $subtitle = $candidate
        break
    }
}
if ($subtitle -eq $null) {
    $subtitle = $d.Paragraphs.Item(2)
}

# Insert a new "Author" paragraph right after the subtitle for the first author.
$subtitle.Range.InsertParagraphAfter()
$author1 = $subtitle.Next()
$author1.Style = "Author"
$author1.Range.Text = "Ben Jarman"

# Insert a second new "Author" paragraph right after the first author.
$author1.Range.InsertParagraphAfter()
$author2 = $author1.Next()
$author2.Style = "Author"
$author2.Range.Text = "Helen Fair"
